$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update therapy activity counts (weekly figures as of third week of September '23)
$ws.Range("C3").Value = 160
$ws.Range("C4").Value = 364
$ws.Range("C5").Value = 131
